$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 79265.7
$ws.Range("I106").Value = 1653.3334
$ws.Range("K106").Value = 1653.3334
$ws.Range("M106").Value = -1022.3334
$ws.Range("H112").Value = 1110
$ws.Range("J112").Value = 1151.25
$ws.Range("L112").Value = 3453.75
$ws.Range("N112").Value = -5669.75
$ws.Range("H113").Value = 1992.2307
$ws.Range("I113").Value = 2255.4443
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 2255.4443
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 998.5556999999999
$ws.Range("N113").Value = -7908
$ws.Range("H134").Value = 56424.875
$ws.Range("J134").Value = 56424.875
$ws.Range("L134").Value = 56424.875
$ws.Range("N134").Value = -66564.875
$ws.Range("H137").Value = 3074.2073
$ws.Range("I137").Value = 1172.2858
$ws.Range("K137").Value = 3516.8574
$ws.Range("M137").Value = -966.8574000000003
$ws.Range("H138").Value = 3037.641
$ws.Range("I138").Value = 2617.913
$ws.Range("J138").Value = 3213.1636
$ws.Range("K138").Value = 7853.739
$ws.Range("L138").Value = 9639.4908
$ws.Range("M138").Value = -2713.739
$ws.Range("N138").Value = -19919.4908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1726.4166
$ws.Range("I45").Value = 1521
$ws.Range("J45").Value = 2753.5
$ws.Range("K45").Value = 1521
$ws.Range("L45").Value = 2753.5
$ws.Range("M45").Value = -1144
$ws.Range("N45").Value = -3507.5
$ws.Range("H61").Value = 2172.55
$ws.Range("I61").Value = 1998.2084
$ws.Range("J61").Value = 2434.0625
$ws.Range("K61").Value = 1998.2084
$ws.Range("L61").Value = 2434.0625
$ws.Range("M61").Value = -1786.2084
$ws.Range("N61").Value = -2858.0625
$ws.Range("H122").Value = 1628
$ws.Range("I122").Value = 1504
$ws.Range("K122").Value = 4512
$ws.Range("M122").Value = -2062
$ws.Range("H131").Value = 50037.668
$ws.Range("J131").Value = 50037.668
$ws.Range("L131").Value = 50037.668
$ws.Range("N131").Value = -60117.668
$ws.Range("H132").Value = 9261435
$ws.Range("I132").Value = 15152994
$ws.Range("J132").Value = 3271.1428
$ws.Range("K132").Value = 45458982
$ws.Range("L132").Value = 9813.428400000001
$ws.Range("M132").Value = -45456452
$ws.Range("N132").Value = -14873.4284
$ws.Range("H136").Value = 2172.55
$ws.Range("I136").Value = 1998.2084
$ws.Range("J136").Value = 2434.0625
$ws.Range("K136").Value = 5994.6252
$ws.Range("L136").Value = 7302.1875
$ws.Range("M136").Value = -3444.6252
$ws.Range("N136").Value = -12402.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1256.8
$ws.Range("I107").Value = 912.5
$ws.Range("J107").Value = 1945.4
$ws.Range("K107").Value = 912.5
$ws.Range("L107").Value = 1945.4
$ws.Range("M107").Value = 1007.5
$ws.Range("N107").Value = -5785.4
$ws.Range("H130").Value = 47019.75
$ws.Range("J130").Value = 47019.75
$ws.Range("L130").Value = 47019.75
$ws.Range("N130").Value = -57059.75
$ws.Range("H135").Value = 52616.5
$ws.Range("J135").Value = 52616.5
$ws.Range("L135").Value = 52616.5
$ws.Range("N135").Value = -62756.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9954.154
$ws.Range("I31").Value = 2169.1
$ws.Range("J31").Value = 35904.332
$ws.Range("K31").Value = 2169.1
$ws.Range("L31").Value = 35904.332
$ws.Range("M31").Value = -1874.1
$ws.Range("N31").Value = -36494.332
$ws.Range("H34").Value = 9954.154
$ws.Range("I34").Value = 2169.1
$ws.Range("J34").Value = 35904.332
$ws.Range("K34").Value = 2169.1
$ws.Range("L34").Value = 35904.332
$ws.Range("M34").Value = -1967.1
$ws.Range("N34").Value = -36308.332
$ws.Range("H99").Value = 2446.875
$ws.Range("I99").Value = 2440.9092
$ws.Range("J99").Value = 2460
$ws.Range("K99").Value = 2440.9092
$ws.Range("L99").Value = 2460
$ws.Range("M99").Value = -942.9092000000001
$ws.Range("N99").Value = -5456
$ws.Range("H111").Value = 47672.5
$ws.Range("J111").Value = 47672.5
$ws.Range("L111").Value = 47672.5
$ws.Range("N111").Value = -55852.5
$ws.Range("H122").Value = 166870600
$ws.Range("I122").Value = 250300660
$ws.Range("K122").Value = 750901980
$ws.Range("M122").Value = -750899530
$ws.Range("H126").Value = 2446.875
$ws.Range("I126").Value = 2440.9092
$ws.Range("J126").Value = 2460
$ws.Range("K126").Value = 7322.7276
$ws.Range("L126").Value = 7380
$ws.Range("M126").Value = -4852.7276
$ws.Range("N126").Value = -12320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2287.6667
$ws.Range("I64").Value = 1803
$ws.Range("J64").Value = 3257
$ws.Range("K64").Value = 5409
$ws.Range("L64").Value = 9771
$ws.Range("M64").Value = -5139
$ws.Range("N64").Value = -10311
$ws.Range("H67").Value = 2287.6667
$ws.Range("I67").Value = 1803
$ws.Range("J67").Value = 3257
$ws.Range("K67").Value = 5409
$ws.Range("L67").Value = 9771
$ws.Range("M67").Value = -4473
$ws.Range("N67").Value = -11643
$ws.Range("H95").Value = 901841.3
$ws.Range("I95").Value = 524
$ws.Range("J95").Value = 1352500
$ws.Range("K95").Value = 1572
$ws.Range("L95").Value = 4057500
$ws.Range("M95").Value = 487
$ws.Range("N95").Value = -4061618
$ws.Range("H131").Value = 46971.72
$ws.Range("I131").Value = 17271.5
$ws.Range("J131").Value = 51426.75
$ws.Range("K131").Value = 51814.5
$ws.Range("L131").Value = 154280.25
$ws.Range("M131").Value = -46774.5
$ws.Range("N131").Value = -164360.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 389038.47
$ws.Range("I80").Value = 719714.3
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 719714.3
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -718716.3
$ws.Range("N80").Value = -5246
$ws.Range("H83").Value = 389038.47
$ws.Range("I83").Value = 719714.3
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 3598571.5
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -3593579.5
$ws.Range("N83").Value = -26234
$ws.Range("H102").Value = 1927.9231
$ws.Range("I102").Value = 1705.4445
$ws.Range("J102").Value = 2428.5
$ws.Range("K102").Value = 1705.4445
$ws.Range("L102").Value = 2428.5
$ws.Range("M102").Value = -83.44450000000006
$ws.Range("N102").Value = -5672.5
$ws.Range("H122").Value = 1640.7
$ws.Range("I122").Value = 1600.875
$ws.Range("K122").Value = 4802.625
$ws.Range("M122").Value = -2352.625
$ws.Range("H126").Value = 8200.6875
$ws.Range("I126").Value = 11028.272
$ws.Range("J126").Value = 1980
$ws.Range("K126").Value = 33084.81600000001
$ws.Range("L126").Value = 5940
$ws.Range("M126").Value = -30614.81600000001
$ws.Range("N126").Value = -10880
$ws.Range("H130").Value = 50284
$ws.Range("J130").Value = 50284
$ws.Range("L130").Value = 50284
$ws.Range("N130").Value = -60324
$ws.Range("H132").Value = 4717.968
$ws.Range("I132").Value = 4876.421
$ws.Range("J132").Value = 4467.0835
$ws.Range("K132").Value = 14629.263
$ws.Range("L132").Value = 13401.2505
$ws.Range("M132").Value = -12099.263
$ws.Range("N132").Value = -18461.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2944.2727
$ws.Range("I7").Value = 2763.6667
$ws.Range("J7").Value = 3161
$ws.Range("K7").Value = 2763.6667
$ws.Range("L7").Value = 3161
$ws.Range("M7").Value = -2651.6667
$ws.Range("N7").Value = -3385
$ws.Range("H40").Value = 2543
$ws.Range("I40").Value = 2440.6365
$ws.Range("K40").Value = 2440.6365
$ws.Range("M40").Value = -2304.6365
$ws.Range("H122").Value = 28587.764
$ws.Range("I122").Value = 34204.355
$ws.Range("K122").Value = 102613.065
$ws.Range("M122").Value = -100163.065
$ws.Range("H126").Value = 2944.2727
$ws.Range("I126").Value = 2763.6667
$ws.Range("J126").Value = 3161
$ws.Range("K126").Value = 8291.000100000001
$ws.Range("L126").Value = 9483
$ws.Range("M126").Value = -5821.000100000001
$ws.Range("N126").Value = -14423

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 105715180
$ws.Range("I122").Value = 132143730
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 396431190
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -396428740
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 5883313
$ws.Range("I126").Value = 14706282
$ws.Range("K126").Value = 44118846
$ws.Range("M126").Value = -44116376
$ws.Range("H131").Value = 46601.2
$ws.Range("J131").Value = 46601.2
$ws.Range("L131").Value = 46601.2
$ws.Range("N131").Value = -56681.2
$ws.Range("H135").Value = 64266.332
$ws.Range("J135").Value = 64266.332
$ws.Range("L135").Value = 64266.332
$ws.Range("N135").Value = -74406.33199999999
$ws.Range("H136").Value = 17540.273
$ws.Range("I136").Value = 26690.82
$ws.Range("J136").Value = 2024.1305
$ws.Range("K136").Value = 80072.45999999999
$ws.Range("L136").Value = 6072.3915
$ws.Range("M136").Value = -77522.45999999999
$ws.Range("N136").Value = -11172.3915
